# 05项目计划表.xlsx - add the "第八周周四" (2018.10.25) weekly block and
# update the completion status / notes of the previous ("第八周周三") block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Part A: update rows 84-88 (the 2018.10.24 block) - completion status
#          and extra notes, plus two cells of richer progress text.
# ---------------------------------------------------------------------

# Row 84 - 邱志鹏: completion status becomes a two-run rich string "未完成(40%)"
# and a note "完成github与IDE的相关配置" is added in column D.
$ws.Range("C84").Value = "未完成(40%)"
$ws.Range("C84").Characters(5, 4).Font.Name = "宋体"
$ws.Range("C84").Characters(5, 4).Font.Size = 11

$ws.Range("D84").Value = "完成github与IDE的相关配置"
$ws.Range("D84").Characters(7, 11).Font.Name = "宋体"
$ws.Range("D84").Characters(7, 11).Font.Size = 11

# Row 85 - 黄立根: completion status "未完成（60%）"
$ws.Range("C85").Value = "未完成（60%）"

# Row 86 - 黄俊贤: plan content gains a "demo" suffix, completion "未完成（60%）"
$ws.Range("B86").Value = "完成百度地图的多点标记、和标记连线demo"
$ws.Range("C86").Value = "未完成（60%）"

# Row 87 - 李达波: completion status "未完成（70%）"
$ws.Range("C87").Value = "未完成（70%）"

# Row 88 - 冯德志: completion status "未完成（40%）"
$ws.Range("C88").Value = "未完成（40%）"

# The C84:C88 cells (and D84) switch from the plain bordered style to the
# "宋体" font style already used by column B in this block (e.g. B85).
$ws.Range("B85").Copy()
$ws.Range("C84,D84,C85,C86,C87,C88").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Part B: append the new weekly block, rows 92-100 (2018.10.25, 第八周周四)
# ---------------------------------------------------------------------

# Pre-populate the block with the same formatting as the block directly
# above it (rows 82-90), then overwrite the text per-cell.
$ws.Range("A82:D90").Copy()
$ws.Range("A92:D100").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A92").Value = "日期：2018.10.25 第八周周四"
$ws.Range("B92").Value = ""
$ws.Range("C92").Value = ""
$ws.Range("D92").Value = ""
$ws.Range("A92:D92").Merge()

$ws.Range("A93").Value = "组员"
$ws.Range("B93").Value = "计划内容"
$ws.Range("C93").Value = "完成情况"
$ws.Range("D93").Value = "备注"

$ws.Range("A94").Value = "邱志鹏"
$ws.Range("B94").Value = "继续完成注册第二个界面、修改密码两个界面"
$ws.Range("C94").Value = "未完成(40%)"
$ws.Range("C94").Characters(5, 4).Font.Name = "宋体"
$ws.Range("C94").Characters(5, 4).Font.Size = 11
$ws.Range("D94").Value = ""

$ws.Range("A95").Value = "黄立根"
$ws.Range("B95").Value = "继续完成pc端显示数据的表单分页"
$ws.Range("C95").Value = "未完成（60%）"
$ws.Range("D95").Value = ""

$ws.Range("A96").Value = "黄俊贤"
$ws.Range("B96").Value = "继续完成百度地图的多点标记、和标记连线demo"
$ws.Range("C96").Value = "未完成（60%）"
$ws.Range("D96").Value = ""

$ws.Range("A97").Value = "李达波"
$ws.Range("B97").Value = "继续昨天未完成的界面，加上个人信息以及聊天信息主界面"
$ws.Range("C97").Value = "未完成（70%）"
$ws.Range("D97").Value = ""

$ws.Range("A98").Value = "冯德志"
$ws.Range("B98").Value = "继续完成地图的路线规划"
$ws.Range("C98").Value = "未完成（40%）"
$ws.Range("D98").Value = ""

$ws.Range("A99").Value = "总结："
$ws.Range("B99").Value = ""
$ws.Range("C99").Value = ""
$ws.Range("D99").Value = ""
$ws.Range("A100").Value = ""
$ws.Range("B100").Value = ""
$ws.Range("C100").Value = ""
$ws.Range("D100").Value = ""
$ws.Range("A99:D100").Merge()

# Row 97's plan-content cell (B97) uses the "宋体" font style like the rest
# of the new block, whereas its donor row (B87) in the block above kept the
# plain style - fix it up explicitly.
$ws.Range("B85").Copy()
$ws.Range("B97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Part C: move the selection to the newly added header row, matching
#          where Excel would land right after typing the new block.
# ---------------------------------------------------------------------
$ws.Range("A92:D92").Select()
